# Generate Report for Handback
# - zh-cn file (2b6ee2ab...) handback finished: status "Ready for handoff" ->
#   "Handed back: in sync with en-US", handback datetime refreshed, the old
#   "version mismatch" error cleared.
# - de-de file (2b6ee2ab...) handback finished the same way, at a later
#   timestamp.
# Rows for the other source file (3112559c...) just pick up the refreshed
# "Status" text and (for zh-cn) the shared handback datetime, since they
# shared the same status/datetime strings as row 2 before this edit.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = $newStatus   # E2 (zh-cn)
$wsOverview.Cells.Item(2, 6).Value = $newStatus   # F2 (de-de)
$wsOverview.Cells.Item(3, 5).Value = $newStatus   # E3 (zh-cn)
$wsOverview.Cells.Item(3, 6).Value = $newStatus   # F3 (de-de)

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(2, 3).Value = $newStatus              # C2 Status
$wsZh.Cells.Item(3, 3).Value = $newStatus              # C3 Status
$wsZh.Cells.Item(2, 11).Value = "2016-08-31 08:33:19"  # K2 Latest Handback DateTime
$wsZh.Cells.Item(3, 11).Value = "2016-08-31 08:33:19"  # K3 Latest Handback DateTime
$wsZh.Cells.Item(2, 16).Value = ""                     # P2 Error Detail cleared

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(2, 3).Value = $newStatus              # C2 Status
$wsDe.Cells.Item(3, 3).Value = $newStatus              # C3 Status
$wsDe.Cells.Item(2, 11).Value = "2016-08-31 08:33:38"  # K2 Latest Handback DateTime
$wsDe.Cells.Item(3, 11).Value = "2016-08-31 08:33:38"  # K3 Latest Handback DateTime
$wsDe.Cells.Item(2, 16).Value = ""                     # P2 Error Detail cleared

# ---- Autofit columns affected by the status/error-detail text changes ----
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsZh.Columns.Item(16).AutoFit() | Out-Null
$wsDe.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(16).AutoFit() | Out-Null
